$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 3342.6667
$ws.Range("I19").Value = 1793.2
$ws.Range("J19").Value = 4117.4
$ws.Range("K19").Value = 1793.2
$ws.Range("L19").Value = 4117.4
$ws.Range("M19").Value = -1618.2
$ws.Range("N19").Value = -4467.4
$ws.Range("H20").Value = 495.5
$ws.Range("I20").Value = 495.5
$ws.Range("K20").Value = 495.5
$ws.Range("M20").Value = -265.5
$ws.Range("H35").Value = 495.5
$ws.Range("I35").Value = 495.5
$ws.Range("K35").Value = 495.5
$ws.Range("M35").Value = -116.5
$ws.Range("H53").Value = 943.6667
$ws.Range("I53").Value = 731.5
$ws.Range("J53").Value = 1686.25
$ws.Range("K53").Value = 731.5
$ws.Range("L53").Value = 1686.25
$ws.Range("M53").Value = -94.5
$ws.Range("N53").Value = -2960.25
$ws.Range("H80").Value = 2976.6667
$ws.Range("I80").Value = 2406.5
$ws.Range("J80").Value = 3628.2856
$ws.Range("K80").Value = 7219.5
$ws.Range("L80").Value = 10884.8568
$ws.Range("M80").Value = -6221.5
$ws.Range("N80").Value = -12880.8568
$ws.Range("H83").Value = 2976.6667
$ws.Range("I83").Value = 2406.5
$ws.Range("J83").Value = 3628.2856
$ws.Range("K83").Value = 21658.5
$ws.Range("L83").Value = 32654.5704
$ws.Range("M83").Value = -16666.5
$ws.Range("N83").Value = -42638.5704
$ws.Range("H98").Value = 1374.5
$ws.Range("I98").Value = 1374.5
$ws.Range("K98").Value = 1374.5
$ws.Range("M98").Value = 123.5
$ws.Range("H111").Value = 2637.6667
$ws.Range("I111").Value = 2457.8635
$ws.Range("K111").Value = 7373.5905
$ws.Range("M111").Value = -4306.5905
$ws.Range("H122").Value = 1374.5
$ws.Range("I122").Value = 1374.5
$ws.Range("K122").Value = 4123.5
$ws.Range("M122").Value = -1673.5
$ws.Range("H131").Value = 5626.1055
$ws.Range("I131").Value = 3699
$ws.Range("J131").Value = 7360.5
$ws.Range("K131").Value = 11097
$ws.Range("L131").Value = 22081.5
$ws.Range("M131").Value = -6057
$ws.Range("N131").Value = -32161.5
$ws.Range("H138").Value = 2452.3438
$ws.Range("I138").Value = 1478.3334
$ws.Range("K138").Value = 4435.0002
$ws.Range("M138").Value = 704.9997999999996

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7162.7075
$ws.Range("I32").Value = 6411.316
$ws.Range("K32").Value = 6411.316
$ws.Range("M32").Value = -6124.316
$ws.Range("H97").Value = 1208.12
$ws.Range("I97").Value = 965.2381
$ws.Range("K97").Value = 965.2381
$ws.Range("M97").Value = -469.2381
$ws.Range("H122").Value = 2627.3333
$ws.Range("I122").Value = 2330.3635
$ws.Range("K122").Value = 6991.0905
$ws.Range("M122").Value = -4541.0905

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 369.35294
$ws.Range("I94").Value = 330
$ws.Range("K94").Value = 330
$ws.Range("M94").Value = 121
$ws.Range("H105").Value = 6495.6665
$ws.Range("I105").Value = 7815.857
$ws.Range("K105").Value = 7815.857
$ws.Range("M105").Value = -6068.857

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H5").Value = 1131.1666
$ws.Range("I5").Value = 1601.75
$ws.Range("J5").Value = 190
$ws.Range("K5").Value = 1601.75
$ws.Range("L5").Value = 190
$ws.Range("M5").Value = -1489.75
$ws.Range("N5").Value = -414
$ws.Range("H6").Value = 2999
$ws.Range("I6").Value = 2999
$ws.Range("K6").Value = 2999
$ws.Range("M6").Value = -2886
$ws.Range("H16").Value = 2328.0833
$ws.Range("I16").Value = 1013.375
$ws.Range("K16").Value = 1013.375
$ws.Range("M16").Value = -726.375
$ws.Range("H19").Value = 734.4
$ws.Range("I19").Value = 749.1429000000001
$ws.Range("K19").Value = 749.1429000000001
$ws.Range("M19").Value = -579.1429000000001
$ws.Range("H24").Value = 734.4
$ws.Range("I24").Value = 749.1429000000001
$ws.Range("K24").Value = 749.1429000000001
$ws.Range("M24").Value = -579.1429000000001
$ws.Range("H31").Value = 6315.067
$ws.Range("I31").Value = 4155.636
$ws.Range("K31").Value = 4155.636
$ws.Range("M31").Value = -3860.636
$ws.Range("H34").Value = 6315.067
$ws.Range("I34").Value = 4155.636
$ws.Range("K34").Value = 4155.636
$ws.Range("M34").Value = -3953.636
$ws.Range("H105").Value = 25433.572
$ws.Range("I105").Value = 14602.6
$ws.Range("K105").Value = 14602.6
$ws.Range("M105").Value = -12855.6
$ws.Range("H113").Value = 2328.0833
$ws.Range("I113").Value = 1013.375
$ws.Range("K113").Value = 1013.375
$ws.Range("M113").Value = 1156.625
$ws.Range("H134").Value = 1726.4615
$ws.Range("I134").Value = 1586.3636
$ws.Range("K134").Value = 4759.0908
$ws.Range("M134").Value = -2224.0908

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2353.353
$ws.Range("H127").Value = 4500
$ws.Range("J127").Value = 4500
$ws.Range("L127").Value = 13500
$ws.Range("N127").Value = -23420
$ws.Range("H132").Value = 6247.8335
$ws.Range("I132").Value = 2800
$ws.Range("J132").Value = 7971.75
$ws.Range("K132").Value = 25200
$ws.Range("L132").Value = 71745.75
$ws.Range("M132").Value = -22670
$ws.Range("N132").Value = -76805.75
$ws.Range("H135").Value = 2353.353
$ws.Range("H138").Value = 4213.273
$ws.Range("I138").Value = 4660
$ws.Range("J138").Value = 3022
$ws.Range("K138").Value = 13980
$ws.Range("L138").Value = 9066
$ws.Range("M138").Value = -8840
$ws.Range("N138").Value = -19346

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 13333502
$ws.Range("J2").Value = 499.33334
$ws.Range("L2").Value = 499.33334
$ws.Range("N2").Value = -725.33334
$ws.Range("H102").Value = 2828.6365
$ws.Range("I102").Value = 1634.5
$ws.Range("K102").Value = 1634.5
$ws.Range("M102").Value = -12.5
$ws.Range("H113").Value = 3707.3953
$ws.Range("I113").Value = 3297.2
$ws.Range("J113").Value = 4277.1113
$ws.Range("K113").Value = 3297.2
$ws.Range("L113").Value = 4277.1113
$ws.Range("M113").Value = -1127.2
$ws.Range("N113").Value = -8617.1113
$ws.Range("H122").Value = 3938.6667
$ws.Range("I122").Value = 3558.1667
$ws.Range("K122").Value = 10674.5001
$ws.Range("M122").Value = -8224.500100000001
$ws.Range("H126").Value = 3614.2856
$ws.Range("I126").Value = 600
$ws.Range("K126").Value = 1800
$ws.Range("M126").Value = 670

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H18").Value = 30005
$ws.Range("I18").Value = 30005
$ws.Range("K18").Value = 30005
$ws.Range("M18").Value = -29833
$ws.Range("H40").Value = 3268.5715
$ws.Range("I40").Value = 3338.3333
$ws.Range("J40").Value = 2850
$ws.Range("K40").Value = 3338.3333
$ws.Range("L40").Value = 2850
$ws.Range("M40").Value = -3202.3333
$ws.Range("N40").Value = -3122
$ws.Range("H46").Value = 2144.96
$ws.Range("I46").Value = 732.2222
$ws.Range("K46").Value = 732.2222
$ws.Range("M46").Value = -544.2222
$ws.Range("H61").Value = 2098.5151
$ws.Range("I61").Value = 1472.3684
$ws.Range("K61").Value = 1472.3684
$ws.Range("M61").Value = -1270.3684
$ws.Range("H68").Value = 2488.5557
$ws.Range("I68").Value = 1649.75
$ws.Range("K68").Value = 1649.75
$ws.Range("M68").Value = -900.75
$ws.Range("H71").Value = 2488.5557
$ws.Range("I71").Value = 1649.75
$ws.Range("K71").Value = 8248.75
$ws.Range("M71").Value = -4504.75
$ws.Range("H93").Value = 2952.6365
$ws.Range("J93").Value = 3211.3076
$ws.Range("L93").Value = 3211.3076
$ws.Range("N93").Value = -5707.3076
$ws.Range("H113").Value = 2098.5151
$ws.Range("I113").Value = 1472.3684
$ws.Range("K113").Value = 1472.3684
$ws.Range("M113").Value = 697.6315999999999
$ws.Range("H132").Value = 15613.667
$ws.Range("I132").Value = 14025.105
$ws.Range("J132").Value = 17389.117
$ws.Range("K132").Value = 42075.315
$ws.Range("L132").Value = 52167.351
$ws.Range("M132").Value = -39545.315
$ws.Range("N132").Value = -57227.351

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1513.6666
$ws.Range("I81").Value = 1299.7778
$ws.Range("K81").Value = 2599.5556
$ws.Range("M81").Value = -1538.5556
$ws.Range("H84").Value = 1513.6666
$ws.Range("I84").Value = 1299.7778
$ws.Range("K84").Value = 12997.778
$ws.Range("M84").Value = -7693.778
$ws.Range("H96").Value = 1873.75
$ws.Range("I96").Value = 1831.6666
$ws.Range("J96").Value = 2000
$ws.Range("K96").Value = 1831.6666
$ws.Range("L96").Value = 2000
$ws.Range("M96").Value = -458.6666
$ws.Range("N96").Value = -4746
$ws.Range("H126").Value = 15399.25
$ws.Range("I126").Value = 15399.25
$ws.Range("K126").Value = 46197.75
$ws.Range("M126").Value = -43727.75
$ws.Range("H132").Value = 3070.889
$ws.Range("I132").Value = 3351.4
$ws.Range("J132").Value = 1668.3334
$ws.Range("K132").Value = 10054.2
$ws.Range("L132").Value = 5005.0002
$ws.Range("M132").Value = -7524.200000000001
$ws.Range("N132").Value = -10065.0002
$ws.Range("H133").Value = 60249.168
$ws.Range("J133").Value = 60249.168
$ws.Range("L133").Value = 60249.168
$ws.Range("N133").Value = -70369.16800000001
